# Auto-generated edit script: refresh market-price derived columns
# (currentAveragePrice / NQ / HQ and LeveProfit columns) per the scheduled
# market-data runner. Mirrors diff of Sheets/Ixion_Profits.xlsx.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3602.7715
$ws.Range("I64").Value = 3714.9812
$ws.Range("K64").Value = 3714.9812
$ws.Range("M64").Value = -3466.9812
$ws.Range("H67").Value = 3602.7715
$ws.Range("I67").Value = 3714.9812
$ws.Range("K67").Value = 3714.9812
$ws.Range("M67").Value = -2856.9812
$ws.Range("H98").Value = 1109.1111
$ws.Range("I98").Value = 959.0909
$ws.Range("K98").Value = 959.0909
$ws.Range("M98").Value = 538.9091
$ws.Range("H122").Value = 1109.1111
$ws.Range("I122").Value = 959.0909
$ws.Range("K122").Value = 2877.2727
$ws.Range("M122").Value = -427.2727
$ws.Range("H127").Value = 839.3043
$ws.Range("I127").Value = 537.0526
$ws.Range("J127").Value = 2275
$ws.Range("K127").Value = 1611.1578
$ws.Range("L127").Value = 6825
$ws.Range("M127").Value = 3348.8422
$ws.Range("N127").Value = -16745

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7901.246
$ws.Range("I32").Value = 4965.1895
$ws.Range("K32").Value = 4965.1895
$ws.Range("M32").Value = -4678.1895
$ws.Range("H61").Value = 7493.2856
$ws.Range("I61").Value = 9391.357
$ws.Range("K61").Value = 9391.357
$ws.Range("M61").Value = -9179.357
$ws.Range("H74").Value = 1780.6364
$ws.Range("I74").Value = 1620.2667
$ws.Range("J74").Value = 2124.2856
$ws.Range("K74").Value = 1620.2667
$ws.Range("L74").Value = 2124.2856
$ws.Range("M74").Value = -746.2666999999999
$ws.Range("N74").Value = -3872.2856
$ws.Range("H77").Value = 1780.6364
$ws.Range("I77").Value = 1620.2667
$ws.Range("J77").Value = 2124.2856
$ws.Range("K77").Value = 8101.3335
$ws.Range("L77").Value = 10621.428
$ws.Range("M77").Value = -3733.3335
$ws.Range("N77").Value = -19357.428
$ws.Range("H132").Value = 4317.8
$ws.Range("I132").Value = 2683.2856
$ws.Range("K132").Value = 8049.8568
$ws.Range("M132").Value = -5519.8568
$ws.Range("H136").Value = 7493.2856
$ws.Range("I136").Value = 9391.357
$ws.Range("K136").Value = 28174.071
$ws.Range("M136").Value = -25624.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11233.2
$ws.Range("I20").Value = 1051.75
$ws.Range("J20").Value = 20631.46
$ws.Range("K20").Value = 1051.75
$ws.Range("L20").Value = 20631.46
$ws.Range("M20").Value = -804.75
$ws.Range("N20").Value = -21125.46
$ws.Range("H102").Value = 2556
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H134").Value = 7633.3335
$ws.Range("I134").Value = 10486.385
$ws.Range("K134").Value = 31459.155
$ws.Range("M134").Value = -28924.155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4288.421
$ws.Range("I31").Value = 1595.1163
$ws.Range("J31").Value = 12560.714
$ws.Range("K31").Value = 1595.1163
$ws.Range("L31").Value = 12560.714
$ws.Range("M31").Value = -1300.1163
$ws.Range("N31").Value = -13150.714
$ws.Range("H34").Value = 4288.421
$ws.Range("I34").Value = 1595.1163
$ws.Range("J34").Value = 12560.714
$ws.Range("K34").Value = 1595.1163
$ws.Range("L34").Value = 12560.714
$ws.Range("M34").Value = -1393.1163
$ws.Range("N34").Value = -12964.714
$ws.Range("H99").Value = 1165.625
$ws.Range("I99").Value = 1185.1666
$ws.Range("J99").Value = 1107
$ws.Range("K99").Value = 1185.1666
$ws.Range("L99").Value = 1107
$ws.Range("M99").Value = 312.8334
$ws.Range("N99").Value = -4103
$ws.Range("H122").Value = 5559161
$ws.Range("I122").Value = 9261925
$ws.Range("J122").Value = 5014
$ws.Range("K122").Value = 27785775
$ws.Range("L122").Value = 15042
$ws.Range("M122").Value = -27783325
$ws.Range("N122").Value = -19942
$ws.Range("H126").Value = 1165.625
$ws.Range("I126").Value = 1185.1666
$ws.Range("J126").Value = 1107
$ws.Range("K126").Value = 3555.4998
$ws.Range("L126").Value = 3321
$ws.Range("M126").Value = -1085.4998
$ws.Range("N126").Value = -8261
$ws.Range("H132").Value = 2769
$ws.Range("I132").Value = 2552
$ws.Range("J132").Value = 2899.2
$ws.Range("K132").Value = 7656
$ws.Range("L132").Value = 8697.599999999999
$ws.Range("M132").Value = -5126
$ws.Range("N132").Value = -13757.6
$ws.Range("H134").Value = 3526.9375
$ws.Range("I134").Value = 3671.52
$ws.Range("K134").Value = 11014.56
$ws.Range("M134").Value = -8479.559999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1055.1818
$ws.Range("I131").Value = 845
$ws.Range("J131").Value = 1101.8889
$ws.Range("K131").Value = 2535
$ws.Range("L131").Value = 3305.6667
$ws.Range("M131").Value = 2505
$ws.Range("N131").Value = -13385.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 36708296
$ws.Range("I122").Value = 54094780
$ws.Range("J122").Value = 3487.3333
$ws.Range("K122").Value = 162284340
$ws.Range("L122").Value = 10461.9999
$ws.Range("M122").Value = -162281890
$ws.Range("N122").Value = -15361.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 73321.57000000001
$ws.Range("J7").Value = 1799
$ws.Range("K7").Value = 73321.57000000001
$ws.Range("L7").Value = 1799
$ws.Range("M7").Value = -73209.57000000001
$ws.Range("N7").Value = -2023
$ws.Range("H55").Value = 271.76923
$ws.Range("I55").Value = 214.55556
$ws.Range("K55").Value = 214.55556
$ws.Range("M55").Value = -41.55556000000001
$ws.Range("I126").Value = 73321.57000000001
$ws.Range("J126").Value = 1799
$ws.Range("K126").Value = 219964.71
$ws.Range("L126").Value = 5397
$ws.Range("M126").Value = -217494.71
$ws.Range("N126").Value = -10337

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1139.9412
$ws.Range("I122").Value = 1139.9412
$ws.Range("K122").Value = 3419.8236
$ws.Range("M122").Value = -969.8235999999997
$ws.Range("H126").Value = 867.2857
$ws.Range("I126").Value = 672
$ws.Range("J126").Value = 1583.3334
$ws.Range("K126").Value = 2016
$ws.Range("L126").Value = 4750.0002
$ws.Range("M126").Value = 454
$ws.Range("N126").Value = -9690.0002
$ws.Range("H138").Value = 38975
$ws.Range("J138").Value = 38975
$ws.Range("L138").Value = 38975
$ws.Range("N138").Value = -49255
